$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 and J1, matching the formatting of the existing header row (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for new columns I (I0) and J (IF), keyed by row number
$data = @{
    2  = @(1, 7)
    3  = @(1, 5)
    4  = @(1, 4)
    5  = @(1, 5)
    6  = @(1, 6)
    7  = @(1, 5)
    8  = @(1, 4)
    9  = @(1, 4)
    10 = @(1, 5)
    11 = @(5, 6)
    12 = @(7, 7)
    13 = @(4, 5)
    14 = @(6, 6)
    15 = @(6, 7)
    16 = @(9, 9)
    17 = @(7, 7)
    18 = @(5, 6)
    19 = @(9, 9)
    20 = @(6, 6)
    21 = @(6, 7)
    22 = @(8, 8)
    23 = @(8, 8)
    24 = @(8, 8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
